# Added new case for Soroe 1950-2099 GCM4 2.6:
# insert a new "1950" sheet (tree class dbh/height/age/number table) right
# before the existing "1996" sheet, fill it with the new stand data, and
# restore/update the navigation state (selection) on every sheet, leaving
# "2001" as the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet right before "1996" -------------------------
$sheet1996 = $wb.Worksheets.Item("1996")
$new = $wb.Worksheets.Add($sheet1996)
$new.Name = "1950"
$wsDBH = $wb.Worksheets.Item("DBH_m")
$wsHeight = $wb.Worksheets.Item("tree height")
$wsAge = $wb.Worksheets.Item("age")
$wsTreeNum = $wb.Worksheets.Item("tree number")
$ws1996 = $wb.Worksheets.Item("1996")
$ws2001 = $wb.Worksheets.Item("2001")

# --- 2. Header row (same labels used by the other class tables) ----------
$new.Range("A1").Value = "dbh"
$new.Range("B1").Value = "height"
$new.Range("C1").Value = "age"
$new.Range("D1").Value = "number"

# --- 3. Stand data rows 2-8 ------------------------------------------------
$new.Range("A2").Value = 4.7829999999999995
$new.Range("B2").Value = 2.0890124395277598
$new.Range("C2").Value = 29
$new.Range("D2").Value = 138.039870849458

$new.Range("A3").Value = 3.2348000000000003
$new.Range("B3").Value = 2.7335003179871999
$new.Range("C3").Value = 29
$new.Range("D3").Value = 289.88372878386201

$new.Range("A4").Value = 3.6350000000000002
$new.Range("B4").Value = 5.76902186020373
$new.Range("C4").Value = 29
$new.Range("D4").Value = 131.13787730698499

$new.Range("A5").Value = 4.6647692307692301
$new.Range("B5").Value = 10.390201453270199
$new.Range("C5").Value = 29
$new.Range("D5").Value = 296.78572232633502

$new.Range("A6").Value = 8.5336363636363597
$new.Range("B6").Value = 12.670895230873599
$new.Range("C6").Value = 29
$new.Range("D6").Value = 531.45350277041405

$new.Range("A7").Value = 8.6534999999999993
$new.Range("B7").Value = 11.4752878098558
$new.Range("C7").Value = 29
$new.Range("D7").Value = 117.333890222039

$new.Range("A8").Value = 4.3259999999999996
$new.Range("B8").Value = 21.519251902355801
$new.Range("C8").Value = 29
$new.Range("D8").Value = 20.7059806274187

# --- 4. Weighted-average labels + summary row -----------------------------
$new.Range("A10").Value = "media pesata dbh"
$new.Range("B10").Value = "media pesata height"

$new.Range("A11").Formula = "=((A2*D2)+(A3*D3)+(A4*D4)+(A5*D5)+(A6*D6)+(A7*D7)+(A8*D8))/(SUM(D2:D8))"
$new.Range("B11").Formula = "=((B2*D2)+(B3*D3)+(B4*D4)+(B5*D5)+(B6*D6)+(B7*D7)+(B8*D8))/(SUM(D2:D8))"
$new.Range("C11").Value = 29
$new.Range("D11").Formula = "=SUM(D2:D8)"

# --- 5. Recreate per-sheet navigation/selection state ---------------------
[void]$wsDBH.Activate()
[void]$wsDBH.Range("B14").Select()

[void]$wsHeight.Activate()
[void]$wsHeight.Range("B5").Select()

[void]$wsAge.Activate()
[void]$wsAge.Range("AA2:AA8").Select()

[void]$wsTreeNum.Activate()
[void]$wsTreeNum.Range("AA2:AA8").Select()

[void]$new.Activate()
[void]$new.Range("A11").Select()

[void]$ws1996.Activate()
[void]$ws1996.Range("A13").Select()

[void]$ws2001.Activate()
[void]$ws2001.Range("A2:A8").Select()
